$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3015.125
$ws.Range("I116").Value = 2649.2727
$ws.Range("J116").Value = 3820
$ws.Range("K116").Value = 2649.2727
$ws.Range("L116").Value = 3820
$ws.Range("M116").Value = 792.7273
$ws.Range("N116").Value = -10704

$ws.Range("H125").Value = 2437.5
$ws.Range("I125").Value = 1625
$ws.Range("J125").Value = 3250
$ws.Range("K125").Value = 14625
$ws.Range("L125").Value = 29250
$ws.Range("M125").Value = -12165
$ws.Range("N125").Value = -34170

$ws.Range("H138").Value = 2984.08
$ws.Range("I138").Value = 1998.1052
$ws.Range("J138").Value = 3588.3872
$ws.Range("K138").Value = 5994.3156
$ws.Range("L138").Value = 10765.1616
$ws.Range("M138").Value = -854.3155999999999
$ws.Range("N138").Value = -21045.1616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23221.012
$ws.Range("I32").Value = 4203.507
$ws.Range("J32").Value = 102647.06
$ws.Range("K32").Value = 4203.507
$ws.Range("L32").Value = 102647.06
$ws.Range("M32").Value = -3916.507
$ws.Range("N32").Value = -103221.06

$ws.Range("H56").Value = 25055
$ws.Range("I56").Value = 15000
$ws.Range("J56").Value = 35110
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 35110
$ws.Range("M56").Value = -14258
$ws.Range("N56").Value = -36594

$ws.Range("H102").Value = 1263.3334
$ws.Range("I102").Value = 1263.3334
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1263.3334
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 358.6666

$ws.Range("H122").Value = 2504
$ws.Range("I122").Value = 1004.8
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 3014.4
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -564.3999999999996
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1465.3125
$ws.Range("I105").Value = 1442.6333
$ws.Range("J105").Value = 1805.5
$ws.Range("K105").Value = 1442.6333
$ws.Range("L105").Value = 1805.5
$ws.Range("M105").Value = 304.3667
$ws.Range("N105").Value = -5299.5

$ws.Range("H134").Value = 1399.8572
$ws.Range("I134").Value = 1258.8823
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 3776.6469
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -1241.6469
$ws.Range("N134").Value = -11067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 711.8
$ws.Range("I105").Value = 679.7778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 12111.111
$ws.Range("I62").Value = 500
$ws.Range("J62").Value = 13562.5
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 40687.5
$ws.Range("M62").Value = -814
$ws.Range("N62").Value = -42059.5

$ws.Range("H65").Value = 12111.111
$ws.Range("I65").Value = 500
$ws.Range("J65").Value = 13562.5
$ws.Range("K65").Value = 4500
$ws.Range("L65").Value = 122062.5
$ws.Range("M65").Value = -1068
$ws.Range("N65").Value = -128926.5

$ws.Range("H75").Value = 1583.4286
$ws.Range("I75").Value = 1196.8
$ws.Range("J75").Value = 2550
$ws.Range("K75").Value = 3590.4
$ws.Range("L75").Value = 7650
$ws.Range("M75").Value = -2592.4
$ws.Range("N75").Value = -9646

$ws.Range("H78").Value = 1583.4286
$ws.Range("I78").Value = 1196.8
$ws.Range("J78").Value = 2550
$ws.Range("K78").Value = 10771.2
$ws.Range("L78").Value = 22950
$ws.Range("M78").Value = -5779.199999999999
$ws.Range("N78").Value = -32934

$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1877
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3384
$ws.Range("N84").ClearContents()

$ws.Range("H87").Value = 9375.462
$ws.Range("I87").Value = 4716.5454
$ws.Range("J87").Value = 34999.5
$ws.Range("K87").Value = 14149.6362
$ws.Range("L87").Value = 104998.5
$ws.Range("M87").Value = -12901.6362
$ws.Range("N87").Value = -107494.5

$ws.Range("H90").Value = 9375.462
$ws.Range("I90").Value = 4716.5454
$ws.Range("J90").Value = 34999.5
$ws.Range("K90").Value = 42448.9086
$ws.Range("L90").Value = 314995.5
$ws.Range("M90").Value = -36208.9086
$ws.Range("N90").Value = -327475.5

$ws.Range("H92").Value = 541.1818
$ws.Range("I92").Value = 502
$ws.Range("J92").Value = 563.5714
$ws.Range("K92").Value = 1506
$ws.Range("L92").Value = 1690.7142
$ws.Range("M92").Value = -258
$ws.Range("N92").Value = -4186.7142

$ws.Range("H97").Value = 1086.6666
$ws.Range("I97").Value = 240
$ws.Range("J97").Value = 1328.5714
$ws.Range("K97").Value = 720
$ws.Range("L97").Value = 3985.7142
$ws.Range("M97").Value = -224
$ws.Range("N97").Value = -4977.7142

$ws.Range("H98").Value = 435878.75
$ws.Range("I98").Value = 1200.4
$ws.Range("J98").Value = 770246.7
$ws.Range("K98").Value = 3601.2
$ws.Range("L98").Value = 2310740.1
$ws.Range("M98").Value = -2103.2
$ws.Range("N98").Value = -2313736.1

$ws.Range("H104").Value = 2300
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2300
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 6900
$ws.Range("N104").Value = -12142

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

$ws.Range("H107").Value = 596.6875
$ws.Range("I107").Value = 440.5
$ws.Range("J107").Value = 690.4
$ws.Range("K107").Value = 1321.5
$ws.Range("L107").Value = 2071.2
$ws.Range("M107").Value = 598.5
$ws.Range("N107").Value = -5911.2

$ws.Range("H108").Value = 400
$ws.Range("I108").Value = 400
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1200
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1680

$ws.Range("H110").Value = 2885.6667
$ws.Range("I110").Value = 1847.5555
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 5542.666499999999
$ws.Range("L110").Value = 18000
$ws.Range("M110").Value = -1452.666499999999
$ws.Range("N110").Value = -26180

$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 1567

$ws.Range("H117").Value = 5055
$ws.Range("I117").Value = 1539.3334
$ws.Range("J117").Value = 6561.7144
$ws.Range("K117").Value = 4618.0002
$ws.Range("L117").Value = 19685.1432
$ws.Range("M117").Value = -1176.0002
$ws.Range("N117").Value = -26569.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3132.6924
$ws.Range("I122").Value = 1126.1666
$ws.Range("J122").Value = 4852.5713
$ws.Range("K122").Value = 3378.4998
$ws.Range("L122").Value = 14557.7139
$ws.Range("M122").Value = -928.4998000000001
$ws.Range("N122").Value = -19457.7139

$ws.Range("H132").Value = 1966
$ws.Range("I132").Value = 1743.8372
$ws.Range("J132").Value = 2700.8462
$ws.Range("K132").Value = 5231.5116
$ws.Range("L132").Value = 8102.5386
$ws.Range("M132").Value = -2701.5116
$ws.Range("N132").Value = -13162.5386

$ws.Range("H139").Value = 50326
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50326
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50326
$ws.Range("N139").Value = -60606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2383.3333
$ws.Range("I100").Value = 1450
$ws.Range("J100").Value = 2468.182
$ws.Range("K100").Value = 1450
$ws.Range("L100").Value = 2468.182
$ws.Range("M100").Value = -909
$ws.Range("N100").Value = -3550.182

$ws.Range("H132").Value = 3577.9146
$ws.Range("I132").Value = 2228.4807
$ws.Range("J132").Value = 5916.933
$ws.Range("K132").Value = 6685.4421
$ws.Range("L132").Value = 17750.799
$ws.Range("M132").Value = -4155.4421
$ws.Range("N132").Value = -22810.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2250
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6090
